$d = $word.ActiveDocument

# 1. Update the date/weekday heading.
[void]$d.Content.Find.Execute("2025-10-29 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-10-30 Thursday", 2)

# 2. Update every arithmetic-fact cell in the practice table, in document
#    order (row-major: row 1 col 1..5, row 2 col 1..5, ...). Several source
#    values repeat with different replacements, so a positional walk over
#    the table cells is used instead of a global find/replace.
$newValues = @(
    "99-89=10", "63+0=63", "51-50=1", "10+17=27", "86-5=81",
    "8+5=13", "43-33=10", "72-8=64", "42+33=75", "76-14=62",
    "62-62=0", "63-8=55", "85-1=84", "75-33=42", "59+8=67",
    "89-34=55", "2+14=16", "28+1=29", "21-17=4", "51-4=47",
    "75-44=31", "89-30=59", "63+23=86", "50-23=27", "39+43=82",
    "70+15=85", "56-43=13", "64-39=25", "56-12=44", "65-63=2",
    "79-6=73", "33+20=53", "5+38=43", "97-93=4", "76-37=39",
    "90-23=67", "71-26=45", "35-10=25", "0+2=2", "96-26=70",
    "94-62=32", "83+7=90", "52-1=51", "20+8=28", "26+53=79",
    "76-23=53", "15+67=82", "69+16=85", "26+28=54", "62+20=82",
    "18+58=76", "56-14=42", "47-41=6", "97-67=30", "64-3=61",
    "54+36=90", "63+6=69", "76-17=59", "85-64=21", "43+8=51",
    "36-1=35", "66-49=17", "18+66=84", "35+53=88", "24+59=83",
    "28+61=89", "10+4=14", "71-45=26", "31-17=14", "43-6=37",
    "2+8=10", "74-10=64", "8+21=29", "72-4=68", "99-71=28",
    "67-20=47", "93-72=21", "99-81=18", "97-2=95", "25-19=6",
    "9+47=56", "54+17=71", "77+17=94", "42-35=7", "70+13=83",
    "27-15=12", "79+6=85", "39+23=62", "42+49=91", "33+34=67",
    "20+34=54", "84-36=48", "5+50=55", "63-3=60", "80-78=2",
    "45-34=11", "84-1=83", "72+4=76", "97-8=89", "93-87=6"
)

$t = $d.Tables.Item(1)
$rows = $t.Rows.Count
$cols = $t.Columns.Count
$i = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$i]
        $i++
    }
}

Write-Output ("Updated cells: " + $i)
